# This workbook is a weekly/daily price log. The commit adds one new
# observation (sampled between the existing 44943 / 44754 rows), which
# pushes every following record down by one row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 48, shifting existing rows 48-156 down to 49-157.
$ws.Rows.Item(48).Insert()

# Populate the newly inserted row 48 with the new record.
$ws.Range("A48").Value = 11
$ws.Range("B48").Value = "Vega Monumental Concepción"
$ws.Range("C48").Value = "Bíobío"
$ws.Range("D48").Value = 44987
$ws.Range("E48").Value = 8
$ws.Range("F48").Value = 100112024
$ws.Range("G48").Value = "Choclo"
$ws.Range("H48").Value = "Choclero"
$ws.Range("I48").Value = "Primera"
$ws.Range("J48").Value = 10000
$ws.Range("K48").Value = 500
$ws.Range("L48").Value = 550
$ws.Range("M48").Value = 525
$ws.Range("N48").Value = "$/unidad"
$ws.Range("O48").Value = "Región Metropolitana"
$ws.Range("P48").Value = 525
$ws.Range("Q48").Value = 1
$ws.Range("R48").Value = "Hortaliza"
